# "Spare from BL11I now working"
# Adds a new "Sheet2" (placed after "Sheet1") holding the BL11I spare
# GAP/CEN calibration table plus the fitted polynomial coefficients
# (C1..C6), highlights the chosen working point (row 18) in bold, and
# leaves that new sheet active/selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet goes right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)

# --- Title -----------------------------------------------------------
$ws2.Range("A1").Value = "BL11I Spare"

# --- Table header ------------------------------------------------------
$ws2.Range("A4").Value = "GAP"
$ws2.Range("B4").Value = "CEN"

# --- GAP / CEN data table (rows 7-25) --------------------------------
$data = @(
  @(21.7, 11.5),
  @(20.7, 12),
  @(16.9, 14),
  @(13.4, 16),
  @(10.2, 18),
  @(7.5,  20),
  @(5.3,  22),
  @(3.65, 24),
  @(2.65, 26),
  @(2.25, 28),
  @(2.5,  30),
  @(3.35, 32),
  @(4.75, 34),
  @(6.7,  36),
  @(9.1,  38),
  @(11.9, 40),
  @(15,   42),
  @(18.5, 44),
  @(20.3, 45)
)

$r = 7
foreach ($pair in $data) {
  $ws2.Cells.Item($r, 1).Value = $pair[0]
  $ws2.Cells.Item($r, 2).Value = $pair[1]
  $r = $r + 1
}

# --- Fitted polynomial coefficients (column E, rows 8-13) ------------
$ws2.Range("E8").Value  = "C1=37.5753"
$ws2.Range("E9").Value  = "C2=-0.08591837"
$ws2.Range("E10").Value = "C3=-0.1905368"
$ws2.Range("E11").Value = "C4=0.007939368"
$ws2.Range("E12").Value = "C5=-0.0001079176"
$ws2.Range("E13").Value = "C6=0.0000004957167"

# Scratch value (C6 restated as a plain number) next to the table.
$ws2.Range("D23").Value = 0.0000004957167

# Highlight the chosen working point.
$ws2.Range("A18:B18").Font.Bold = $true

# Make the new sheet the active one, selected at E17.
[void]$ws2.Activate()
[void]$ws2.Range("E17").Select()
